# Update contact data: fix "Roshan " (trailing space) to "Roshan" in the
# "data" worksheet, and move the active selection to D9.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A3").Value = "Roshan"

$ws.Range("D9").Select()
